$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "22.556.58"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.33%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.577.56"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("E5").Value = "  +0.03%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "288.78"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.93%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.3695"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.18%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "48.70"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -2.37%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.3344"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("E10").Value = "  +0.22%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.07479"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -0.89%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "21.02"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -1.24%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.010"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -0.33%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "6.965"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.09%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "1.577.26"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.01%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.00001118"
$c.Style = "Normal"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("E18").Value = "  -2.18%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06762"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.06%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "6.430"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +2.00%  "
$ws.Range("E21").Value = "  +0.04%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "16.59"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.24%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "22.555.46"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +0.34%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.401"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.63%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.604"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.05%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "152.57"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +2.33%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "19.70"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.82%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "5.019"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -0.69%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "124.58"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -0.57%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "1.753.75"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +0.27%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "1.070"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -0.20%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "6.189"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.10%  "
$ws.Range("E34").Value = "  -0.39%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "9.692"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -0.71%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.08317"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.58%  "
$ws.Range("E37").Value = "  -1.31%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.2270"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -1.28%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "5.447"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.06%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.304"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -4.39%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.06395"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -2.36%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "11.41"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.41%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.6364"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +2.16%  "
$ws.Range("E44").Value = "  +0.06%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "14.07"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.44%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.6201"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +5.79%  "
$ws.Range("E47").Value = "  -0.94%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "2.064"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -0.46%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "124.98"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -3.33%  "
$ws.Range("E50").Value = "  -0.19%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.07275"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -0.78%  "
